$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 323 (old rows 323-348 shift down to 324-349)
$ws.Rows(323).Insert()

# Populate the newly inserted row 323 with the new record
$ws.Cells.Item(323, 1).Value = 10
$ws.Cells.Item(323, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(323, 3).Value = "La Araucanía"
$ws.Cells.Item(323, 4).Value = 44461
$ws.Cells.Item(323, 5).Value = 9
$ws.Cells.Item(323, 6).Value = 100112045
$ws.Cells.Item(323, 7).Value = "Zapallo"
$ws.Cells.Item(323, 8).Value = "Paine"
$ws.Cells.Item(323, 9).Value = "1a (guarda)"
$ws.Cells.Item(323, 10).Value = 400
$ws.Cells.Item(323, 11).Value = 500
$ws.Cells.Item(323, 12).Value = 600
$ws.Cells.Item(323, 13).Value = 550
$ws.Cells.Item(323, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(323, 15).Value = "Región del Maule"
$ws.Cells.Item(323, 16).Value = 550
$ws.Cells.Item(323, 17).Value = 1
$ws.Cells.Item(323, 18).Value = "Hortaliza"
